$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# --- Swap Murcia/Pontevedra row order and update their stats ---
# Row 40 was Pontevedra, row 41 was Murcia. After the edit, row 40 becomes
# Murcia (with updated stats) and row 41 becomes Pontevedra (keeping the
# previous Pontevedra stats).
$ws.Range("A40").Value = "Murcia"
$ws.Range("B40").Value = 1579
$ws.Range("C40").Value = 2180
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 149

$ws.Range("A41").Value = "Pontevedra"
$ws.Range("B41").Value = 1536
$ws.Range("C41").Value = 333
$ws.Range("D41").Value = 1411
$ws.Range("E41").Value = 30

# --- Update numeric stats for several other provinces ---
$ws.Range("B4").Value = 67871
$ws.Range("D4").Value = 18158
$ws.Range("E4").Value = 8977

$ws.Range("B5").Value = 57148
$ws.Range("D5").Value = 24244
$ws.Range("E5").Value = 6701

$ws.Range("B6").Value = 18789
$ws.Range("D6").Value = 8112

$ws.Range("B7").Value = 16889
$ws.Range("D7").Value = 7557
$ws.Range("E7").Value = 2940

$ws.Range("B9").Value = 12600
$ws.Range("D9").Value = 538
$ws.Range("E9").Value = 1391

$ws.Range("B13").Value = 5646
$ws.Range("D13").Value = 1016

$ws.Range("B16").Value = 5219
$ws.Range("D16").Value = 805

$ws.Range("D20").Value = 571
$ws.Range("E20").Value = 358

$ws.Range("B32").Value = 2397
$ws.Range("D32").Value = 1026
$ws.Range("E32").Value = 308

$ws.Range("B33").Value = 2322
$ws.Range("D33").Value = 628

$ws.Range("B59").Value = 121
$ws.Range("C59").Value = 125
$ws.Range("D59").Value = 0

# --- Update the "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 18:05"
